$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Varun Chakravarthy"

$ws.Range("A1:M4").NumberFormat = "@"

$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$data = @(
    @("10th","Kolkata Knight Riders","Varun Chakravarthy","","2","3","0","0","66.66","Royal Challengers Bangalore","Chennai","April 18","RCB won by 38 runs"),
    @("Final","Kolkata Knight Riders","Varun Chakravarthy","","0","0","0","0","-","Chennai Super Kings","Dubai (DSC)","October 15","Super Kings won by 27 runs"),
    @("15th","Kolkata Knight Riders","Varun Chakravarthy","run out (Chahar/Curran)","0","1","0","0","0.00","Chennai Super Kings","Wankhede","April 21","Super Kings won by 18 runs")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws.Range("A1:M4").ClearFormats()
